$d = $word.ActiveDocument

# The SOP's author-title line used to read "...Title: Prototyping Labs
# Manager..." in the signature block table, and the document's first-page
# header repeated the same "Prototyping Labs at GIX" wording. Both were
# corrected to the singular "Prototyping Lab" (department name fix).

# 1) Signature-block table cell in the document body (word/document.xml)
$bodyRange = $d.Content
$bodyRange.Find.Execute("Prototyping Labs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Prototyping Lab", 2) | Out-Null

# 2) First-page header title (word/header1.xml)
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        $hdrRange = $hdr.Range
        if ($hdrRange.Text -like "*Prototyping Labs*") {
            $hdrRange.Find.Execute("Prototyping Labs", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "Prototyping Lab", 2) | Out-Null
        }
    }
}
